# Ironed out a few bugs, created 'visible' column.
#
# Inserts a new "Visible" boolean column between the existing "Enabled"
# column (F) and the "Options" column (old G, now H) on the Input sheet,
# filling it with TRUE for every data row, and makes Input the active /
# selected sheet (previously Output was active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Insert a new column at G; everything from the old G onward (Options,
# Errormessage) shifts one column to the right (-> H, I).
$ws.Columns("G:G").Insert()

# New header cell, matching the bold header style used by the rest of row 1.
$ws.Range("G1").Value = "Visible"
$ws.Range("G1").Font.Bold = $true

# New data cells - all rows visible (TRUE) by default.
$ws.Range("G2").Value = $true
$ws.Range("G3").Value = $true
$ws.Range("G4").Value = $true
$ws.Range("G5").Value = $true

# Input becomes the active sheet/tab (was Output before).
$ws.Activate() | Out-Null
$ws.Range("G6").Select() | Out-Null
